$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch a style first so the (originally style-less) workbook has a valid
# cellXfs table before any cell write is attempted.
$ws.Range("A1").Style = "Normal"

# Force column H (operationperformeddatetime) to be stored as plain text,
# matching the original t="str" serial-looking strings.
$ws.Columns.Item(8).NumberFormat = "@"

# Insert two new rows above the current row 6, pushing the existing
# rows 6-10 down to rows 8-12.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# Update the "operationperformeddatetime" column on the rows that kept
# their data but changed date formatting (now serial-like text).
$ws.Range("H2").Value = "42328.6278935185"
$ws.Range("H3").Value = "42328.6278935185"
$ws.Range("H4").Value = "42339.942037037"
$ws.Range("H5").Value = "42326.942037037"

# New row 6: edf6eed1-40f2-469c-ae7e-f6e0de90167a / Episode
# (a lone "'" forces a genuine, present-but-empty text cell instead of no
# cell at all, matching the blank parentid/parents/operationperformedby
# cells that already exist elsewhere in this sheet)
$ws.Range("A6").Value = "update"
$ws.Range("B6").Value = "edf6eed1-40f2-469c-ae7e-f6e0de90167a"
$ws.Range("C6").Value = "Read, Write, Learn"
$ws.Range("D6").Value = "Episode"
$ws.Range("E6").Value = "'"
$ws.Range("F6").Value = '[{"id":"40d2a3eb-40f2-4e10-806c-50f310db5957","subtype":null,"type":"Unit"}]'
$ws.Range("G6").Value = "Learning Analytics"
$ws.Range("H6").Value = "42326.942037037"

# New row 7: 40d2a3eb-40f2-4e10-806c-50f310db5957 / Unit
$ws.Range("A7").Value = "update"
$ws.Range("B7").Value = "40d2a3eb-40f2-4e10-806c-50f310db5957"
$ws.Range("C7").Value = "Read, Write, Learn"
$ws.Range("D7").Value = "Unit"
$ws.Range("E7").Value = "'"
$ws.Range("F7").Value = "'"
$ws.Range("G7").Value = "Learning Analytics"
$ws.Range("H7").Value = "42326.942037037"

# Row 8 (was row 6: 1934288 / Composite / question) -- only the date text changes
$ws.Range("H8").Value = "42328.6303356482"

# Row 9 (was row 7: 1934291 / Composite / question) -- only the date text changes
$ws.Range("H9").Value = "42328.6304166667"

# Row 10 (was row 8: 1934293 / Composite / question) -- only the date text changes
$ws.Range("H10").Value = "42328.6305092593"

# Row 11 (was row 9: c9523d19... / 1.1 Unit Accomplishments / assessment) -- only the date text changes
$ws.Range("H11").Value = "42328.6303356482"

# Row 12 (was row 10: 40d2a3eb... / Read, Write, Learn / Unit) -- only the date text changes
$ws.Range("H12").Value = "42339.942037037"
